$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily profit row (row 74 -> row 75), mirroring the existing
# layout: column A holds the date as plain text, column B the numeric profit.
$row = 75

$cellA = $ws.Cells.Item($row, 1)
# Format as text first so the "10/31/2025" string isn't auto-converted into a
# date serial number, then drop back to the default (Normal) style so the
# cell carries no explicit style index, matching the rest of the data rows.
$cellA.NumberFormat = "@"
$cellA.Value = "10/31/2025"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 10798.86

$wb.Save()
